$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = " 06-11-20"
$ws.Range("A57").Value = '$ 17.323 CLP 06-11-20'
